$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 227; this shifts the existing rows
# 227-302 down to 228-303 (values travel with their rows automatically).
$ws.Rows("227:227").Insert()

# Populate the newly-inserted row 227 with its data. The "constant" columns
# (market/product metadata) match every other row in this block; only the
# measurement columns (D, J, K, L, M, P) plus I/O (which happen to equal the
# old row 227 values) are meaningful here.
$ws.Range("A227").Value = 10
$ws.Range("B227").Value = "Vega Modelo de Temuco"
$ws.Range("C227").Value = "La Araucanía"
$ws.Range("D227").Value = 44900
$ws.Range("E227").Value = 9
$ws.Range("F227").Value = 100112039
$ws.Range("G227").Value = "Ciboulette"
$ws.Range("H227").Value = "Sin especificar"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 55
$ws.Range("K227").Value = 5000
$ws.Range("L227").Value = 5000
$ws.Range("M227").Value = 5000
$ws.Range("N227").Value = "$/docena de atados"
$ws.Range("O227").Value = "Provincia de Cautín"
$ws.Range("P227").Value = 1667
$ws.Range("Q227").Value = 3
$ws.Range("R227").Value = "Hortaliza"

# Match the date-format style used by the rest of column D.
$ws.Range("D227").NumberFormat = "YYYY-MM-DD HH:MM:SS"
